# Apply cell text/value updates to the cryptos worksheet (generated from the OOXML diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain text in the source file (t="inlineStr"), including
# values such as "63.795.64" (dot-grouped, not a valid number) and values like "146.10"
# whose trailing zero would be lost if Excel re-interpreted the text as a number. Force
# each target cell to Text format before writing so the literal string is preserved.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.795.64"
$ws.Range("E2").Value = "  -3.10%  "
$ws.Range("D3").Value = "3.141.24"
$ws.Range("E3").Value = "  -2.94%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "608.27"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").Value = "146.10"
$ws.Range("E6").Value = "  -6.60%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.137.55"
$ws.Range("E8").Value = "  -2.87%  "
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("E10").Value = "  -4.94%  "
$ws.Range("D11").Value = "5.31"
$ws.Range("E11").Value = "  -5.76%  "
$ws.Range("E12").Value = "  -2.93%  "
$ws.Range("D13").Value = "0.0000253"
$ws.Range("E13").Value = "  -4.15%  "
$ws.Range("D14").Value = "35.52"
$ws.Range("E14").Value = "  -6.68%  "
$ws.Range("D15").Value = "3.651.47"
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("D17").Value = "63.831.27"
$ws.Range("E17").Value = "  -3.32%  "
$ws.Range("D18").Value = "3.138.86"
$ws.Range("E18").Value = "  -3.32%  "
$ws.Range("D19").Value = "6.84"
$ws.Range("E19").Value = "  -4.63%  "
$ws.Range("D20").Value = "477.86"
$ws.Range("E20").Value = "  -3.30%  "
$ws.Range("D21").Value = "14.65"
$ws.Range("E21").Value = "  -3.30%  "
$ws.Range("D22").Value = "0.709"
$ws.Range("E22").Value = "  -3.61%  "
$ws.Range("D23").Value = "7.85"
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("D24").Value = "13.57"
$ws.Range("E24").Value = "  -5.36%  "
$ws.Range("D25").Value = "83.77"
$ws.Range("E25").Value = "  -2.92%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").Value = "2.80"
$ws.Range("E27").Value = "  -6.52%  "
$ws.Range("D28").Value = "8.51"
$ws.Range("E28").Value = "  -4.94%  "
$ws.Range("D29").Value = "0.122"
$ws.Range("E29").Value = "  -7.02%  "
$ws.Range("D30").Value = "6.94"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -10.72%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "2.71"
$ws.Range("E33").Value = "  -3.36%  "
$ws.Range("D34").Value = "26.27"
$ws.Range("E34").Value = "  -4.52%  "
$ws.Range("D35").Value = "1.13"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").Value = "0.0₃0796"
$ws.Range("E36").Value = "  +6.12%  "
$ws.Range("D37").Value = "5.96"
$ws.Range("E37").Value = "  -5.42%  "
$ws.Range("D38").Value = "53.07"
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("D39").Value = "459.82"
$ws.Range("E39").Value = "  -5.14%  "
$ws.Range("D40").Value = "3.04"
$ws.Range("E40").Value = "  -8.69%  "
$ws.Range("E41").Value = "  -5.17%  "
$ws.Range("D42").Value = "0.119"
$ws.Range("E42").Value = "  -7.66%  "
$ws.Range("D43").Value = "8.36"
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("D44").Value = "2.848.65"
$ws.Range("E44").Value = "  -4.30%  "
$ws.Range("D45").Value = "2.31"
$ws.Range("E45").Value = "  -7.46%  "
$ws.Range("D46").Value = "0.268"
$ws.Range("E46").Value = "  -5.99%  "
$ws.Range("D47").Value = "2.46"
$ws.Range("E47").Value = "  +1.31%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "26.34"
$ws.Range("E48").Value = "  -6.13%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  -3.41%  "
$ws.Range("D51").Value = "119.02"
$ws.Range("E51").Value = "  -2.37%  "

# Restore the default (unstyled) appearance now that the text values are set, so the
# cells end up without an explicit style index, matching the original workbook.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"

